# Refresh cryptos price/volume table to match the latest scrape.
# Each entry below carries the already-updated cell value; writing through
# a leading apostrophe forces Excel to keep numeric-looking text (e.g. "238.01")
# as a literal string instead of re-interpreting it as a number, and resetting
# the range Style back to "Normal" clears the quote-prefix flag that the
# apostrophe trick leaves behind, so no cell formatting changes either.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '29.244.29' }
    @{ Cell = 'E2'; Value = '  +0.26%  ' }
    @{ Cell = 'D3'; Value = '1.855.58' }
    @{ Cell = 'E3'; Value = '  +0.19%  ' }
    @{ Cell = 'E4'; Value = '  -0.06%  ' }
    @{ Cell = 'D5'; Value = '0.7013' }
    @{ Cell = 'E5'; Value = '  +2.98%  ' }
    @{ Cell = 'D6'; Value = '238.01' }
    @{ Cell = 'E6'; Value = '  +0.18%  ' }
    @{ Cell = 'E7'; Value = '  -0.07%  ' }
    @{ Cell = 'D8'; Value = '0.08052' }
    @{ Cell = 'E8'; Value = '  +4.73%  ' }
    @{ Cell = 'D9'; Value = '0.3025' }
    @{ Cell = 'E9'; Value = '  -0.30%  ' }
    @{ Cell = 'D10'; Value = '23.59' }
    @{ Cell = 'E10'; Value = '  +2.22%  ' }
    @{ Cell = 'D11'; Value = '0.08182' }
    @{ Cell = 'E11'; Value = '  +0.56%  ' }
    @{ Cell = 'B12'; Value = 'Polkadot' }
    @{ Cell = 'C12'; Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot' }
    @{ Cell = 'D12'; Value = '5.205' }
    @{ Cell = 'E12'; Value = '  +0.64%  ' }
    @{ Cell = 'B13'; Value = 'Polygon' }
    @{ Cell = 'C13'; Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic' }
    @{ Cell = 'D13'; Value = '0.7069' }
    @{ Cell = 'E13'; Value = '  -1.99%  ' }
    @{ Cell = 'B14'; Value = 'WrappedEther' }
    @{ Cell = 'C14'; Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth' }
    @{ Cell = 'D14'; Value = '1.775.34' }
    @{ Cell = 'E14'; Value = '  -4.76%  ' }
    @{ Cell = 'D15'; Value = '89.76' }
    @{ Cell = 'E15'; Value = '  +0.68%  ' }
    @{ Cell = 'D16'; Value = '29.276.58' }
    @{ Cell = 'D17'; Value = '5.833' }
    @{ Cell = 'E17'; Value = '  +2.17%  ' }
    @{ Cell = 'D18'; Value = '0.000007862' }
    @{ Cell = 'E18'; Value = '  +0.69%  ' }
    @{ Cell = 'E19'; Value = '  +1.00%  ' }
    @{ Cell = 'D20'; Value = '236.84' }
    @{ Cell = 'E20'; Value = '  +1.61%  ' }
    @{ Cell = 'D22'; Value = '2.105.27' }
    @{ Cell = 'E22'; Value = '  +0.50%  ' }
    @{ Cell = 'E23'; Value = '  -0.03%  ' }
    @{ Cell = 'D24'; Value = '7.487' }
    @{ Cell = 'E24'; Value = '  +0.99%  ' }
    @{ Cell = 'D25'; Value = '163.33' }
    @{ Cell = 'E25'; Value = '  +1.01%  ' }
    @{ Cell = 'D26'; Value = '8.889' }
    @{ Cell = 'E26'; Value = '  -0.61%  ' }
    @{ Cell = 'D27'; Value = '0.1413' }
    @{ Cell = 'E27'; Value = '  -0.87%  ' }
    @{ Cell = 'D28'; Value = '18.09' }
    @{ Cell = 'E28'; Value = '  +0.54%  ' }
    @{ Cell = 'D29'; Value = '1.912' }
    @{ Cell = 'E29'; Value = '  -1.98%  ' }
    @{ Cell = 'D30'; Value = '1.409' }
    @{ Cell = 'E30'; Value = '  +1.44%  ' }
    @{ Cell = 'E31'; Value = '  -0.70%  ' }
    @{ Cell = 'D32'; Value = '4.367' }
    @{ Cell = 'E32'; Value = '  -3.18%  ' }
    @{ Cell = 'D33'; Value = '4.027' }
    @{ Cell = 'E33'; Value = '  +0.57%  ' }
    @{ Cell = 'D34'; Value = '0.05198' }
    @{ Cell = 'E34'; Value = '  +0.84%  ' }
    @{ Cell = 'D35'; Value = '1.165' }
    @{ Cell = 'E35'; Value = '  -1.09%  ' }
    @{ Cell = 'D36'; Value = '0.7164' }
    @{ Cell = 'E36'; Value = '  +2.09%  ' }
    @{ Cell = 'D37'; Value = '0.9971' }
    @{ Cell = 'E37'; Value = '  -2.47%  ' }
    @{ Cell = 'D38'; Value = '2.688' }
    @{ Cell = 'E38'; Value = '  +0.45%  ' }
    @{ Cell = 'D39'; Value = '0.01850' }
    @{ Cell = 'E39'; Value = '  +0.45%  ' }
    @{ Cell = 'E40'; Value = '  +1.79%  ' }
    @{ Cell = 'D41'; Value = '0.9345' }
    @{ Cell = 'E41'; Value = '  +2.76%  ' }
    @{ Cell = 'D42'; Value = '1.152.27' }
    @{ Cell = 'E42'; Value = '  +5.01%  ' }
    @{ Cell = 'D43'; Value = '6.007' }
    @{ Cell = 'E43'; Value = '  +0.35%  ' }
    @{ Cell = 'D44'; Value = '0.4264' }
    @{ Cell = 'E44'; Value = '  +0.00%  ' }
    @{ Cell = 'D45'; Value = '70.16' }
    @{ Cell = 'E45'; Value = '  +0.15%  ' }
    @{ Cell = 'D46'; Value = '1.000' }
    @{ Cell = 'E46'; Value = '  -0.11%  ' }
    @{ Cell = 'D47'; Value = '102.82' }
    @{ Cell = 'E47'; Value = '  +0.70%  ' }
    @{ Cell = 'D48'; Value = '0.5285' }
    @{ Cell = 'E48'; Value = '  -3.37%  ' }
    @{ Cell = 'B49'; Value = 'RenderToken' }
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr' }
    @{ Cell = 'D49'; Value = '1.747' }
    @{ Cell = 'E49'; Value = '  -0.79%  ' }
    @{ Cell = 'B50'; Value = 'RocketPoolETH' }
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth' }
    @{ Cell = 'D50'; Value = '1.998.40' }
    @{ Cell = 'E50'; Value = '  +0.18%  ' }
    @{ Cell = 'D51'; Value = '9.148' }
    @{ Cell = 'E51'; Value = '  +0.26%  ' }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    $range.Value = "'" + $u.Value
    $range.Style = "Normal"
}
